$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the old bold "Play Age of Asgard for Free..." paragraph that used
#    to sit right before the italic "Experience the world..." paragraph at
#    the end of the document, and replace that italic paragraph's text with
#    the new "Prompt: ..." text (keeping its italic formatting intact).
#    Do this FIRST, while the title text is still unique in the document, so
#    the Find/Replace below cannot match anything else.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Play Age of Asgard for Free: Norse Mythology Slot Game`r") {
        $p.Range.Delete()
        break
    }
}

$d.Content.Find.Execute("Experience the world of Norse mythology in Age of Asgard, a unique slot game with two grids and special symbols. Play for free and battle for riches.", $true, $false, $false, $false, $false, $true, 1, $false, "Prompt: Create a feature image for Age of Asgard, a slot game that offers a unique twist on the beloved mythological theme of Norse gods and their battles. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The image should showcase the excitement and epicness of the game, with an ancient Viking village in the background and flames and warriors on both sides of the grid. The Maya warrior should be dressed in a traditional Viking outfit and have a big smile on their face, holding up a winning combination on the slot machine. Be sure to incorporate elements of Norse mythology and the different symbols of the game into the design of the image.", 2)

# ---------------------------------------------------------------------------
# 2. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$insertPoint = $d.Range($metaStart, $metaStart)
$insertPoint.InsertAfter("Meta description: Experience the world of Norse mythology in Age of Asgard, a unique slot game with two grids and special symbols. Play for free and battle for riches.")

$boldRange = $d.Range($metaStart, $metaStart + 16)
$boldRange.Bold = 1
